$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (rich-text run replacements; formatting unchanged) ---
$hdr = $ws.Range("A8")
$hdr.Characters(21, 2).Text = "47"

$week = $ws.Range("C9")
$week.Characters(27, 10).Text = "11/18/2024"
$week.Characters(48, 10).Text = "11/24/2024"

# --- Crime-data table updates (rows 15-33) ---

# Row 15
$ws.Range("D15").NumberFormat = '#,##0'
$ws.Range("D15").Value = 1
$ws.Range("E15").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("E15").Value = -100
$ws.Range("G15").NumberFormat = '#,##0'
$ws.Range("G15").Value = 1
$ws.Range("H15").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 8
$ws.Range("K15").Value = 62.5

# Row 16
$ws.Range("C16").NumberFormat = '#,##0'
$ws.Range("C16").Value = 2
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 100
$ws.Range("I16").Value = 65
$ws.Range("J16").Value = 65
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 32.653061224489
$ws.Range("M16").Value = -30.10752688172
$ws.Range("N16").Value = -87.354085603112

# Row 17
$ws.Range("C17").Value = 3
$ws.Range("D17").Value = 2
$ws.Range("E17").Value = 50
$ws.Range("G17").Value = 9
$ws.Range("H17").Value = 44.444444444444
$ws.Range("I17").Value = 113
$ws.Range("J17").Value = 124
$ws.Range("K17").Value = -8.870967741935
$ws.Range("L17").Value = 17.708333333333
$ws.Range("M17").Value = 20.212765957446
$ws.Range("N17").Value = -56.870229007633

# Row 18
$ws.Range("C18").NumberFormat = '#,##0'
$ws.Range("C18").Value = 1
$ws.Range("F18").Value = 3
$ws.Range("G18").Value = 3
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 81
$ws.Range("K18").Value = -4.705882352941
$ws.Range("L18").Value = -11.95652173913
$ws.Range("M18").Value = -64.159292035398
$ws.Range("N18").Value = -92.011834319526

# Row 19
$ws.Range("C19").Value = 5
$ws.Range("E19").Value = -28.571428571428
$ws.Range("F19").Value = 28
$ws.Range("G19").Value = 27
$ws.Range("H19").Value = 3.703703703703
$ws.Range("I19").Value = 333
$ws.Range("J19").Value = 392
$ws.Range("K19").Value = -15.051020408163
$ws.Range("L19").Value = -28.997867803838
$ws.Range("M19").Value = 13.265306122449
$ws.Range("N19").Value = -21.27659574468

# Row 20
$ws.Range("C20").Value = 1
$ws.Range("D23").Copy($ws.Range("D20"))
$ws.Range("E23").Copy($ws.Range("E20"))
$ws.Range("F20").Value = 7
$ws.Range("G20").Value = 8
$ws.Range("H20").Value = -12.5
$ws.Range("I20").Value = 172
$ws.Range("K20").Value = 47.008547008547
$ws.Range("L20").Value = 79.166666666666
$ws.Range("M20").Value = 24.63768115942
$ws.Range("N20").Value = -90.282485875706

# Row 21
$ws.Range("D21").Value = 11
$ws.Range("E21").Value = 9.090909090909
$ws.Range("F21").Value = 57
$ws.Range("G21").Value = 55
$ws.Range("H21").Value = 3.636363636363
$ws.Range("I21").Value = 778
$ws.Range("J21").Value = 793
$ws.Range("K21").Value = -1.891551071878
$ws.Range("L21").Value = -4.422604422604
$ws.Range("M21").Value = -8.685446009389
$ws.Range("N21").Value = -80.622665006226

# Row 22
$ws.Range("C22").NumberFormat = '#,##0'
$ws.Range("C22").Value = 1
$ws.Range("E22").Value = 0
$ws.Range("F22").NumberFormat = '#,##0'
$ws.Range("F22").Value = 1
$ws.Range("H22").Value = -66.666666666666
$ws.Range("I22").Value = 10
$ws.Range("J22").Value = 10
$ws.Range("L22").Value = -16.666666666666
$ws.Range("M22").Value = 42.857142857142

# Row 24
$ws.Range("C24").Value = 29
$ws.Range("D24").Value = 27
$ws.Range("E24").Value = 7.407407407407
$ws.Range("F24").Value = 108
$ws.Range("G24").Value = 104
$ws.Range("H24").Value = 3.846153846153
$ws.Range("I24").Value = 1284
$ws.Range("J24").Value = 1388
$ws.Range("K24").Value = -7.492795389048
$ws.Range("L24").Value = -24.559341950646
$ws.Range("M24").Value = 17.47483989021

# Row 25
$ws.Range("C25").Value = 15
$ws.Range("D25").Value = 18
$ws.Range("E25").Value = -16.666666666666
$ws.Range("F25").Value = 51
$ws.Range("G25").Value = 65
$ws.Range("H25").Value = -21.538461538461
$ws.Range("I25").Value = 781
$ws.Range("J25").Value = 748
$ws.Range("K25").Value = 4.411764705882
$ws.Range("L25").Value = -24.686595949855

# Row 26
$ws.Range("C26").Value = 8
$ws.Range("D26").Value = 14
$ws.Range("E26").Value = -42.857142857142
$ws.Range("F26").Value = 30
$ws.Range("G26").Value = 31
$ws.Range("H26").Value = -3.225806451612
$ws.Range("I26").Value = 325
$ws.Range("J26").Value = 341
$ws.Range("K26").Value = -4.692082111436
$ws.Range("L26").Value = 14.840989399293
$ws.Range("M26").Value = -8.707865168539

# Row 27
$ws.Range("J27").Value = 14
$ws.Range("K27").Value = 21.428571428571

# Row 28
$ws.Range("D28").NumberFormat = '#,##0'
$ws.Range("D28").Value = 1
$ws.Range("E28").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("E28").Value = 0
$ws.Range("I28").Value = 47
$ws.Range("J28").Value = 27
$ws.Range("K28").Value = 74.074074074074
$ws.Range("L28").Value = 17.5

# Row 33
$ws.Range("D14").Copy($ws.Range("D33"))
$ws.Range("E14").Copy($ws.Range("E33"))

Write-Output "Applied weekly crime-data refresh."
